$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 88 (pushes "SABONETE LIQUIDO..." and everything below down by one)
$ws.Rows.Item(88).Insert()

# Fill in the new product row
$ws.Cells.Item(88, 1).Value = "SABAO EM BARRA - 200G - 200G"
$ws.Cells.Item(88, 2).Value = "UN"
$ws.Cells.Item(88, 3).Value = "S010046"
$ws.Cells.Item(88, 4).Value = 51

# Keep the view similar to the target (scrolled down near the new row)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 85
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H91").Select()

$wb.Save()
